# Weekly update: a new (date, quality) sample pair was inserted at the top
# of the "Betarraga" data block (rows 97-210), pushing every row's
# I/J/K/L/M/O/P (quality/volume/price/origin) values down by one row, while
# the shared Fecha (D) value for each Primera/Segunda pair shifts down by a
# full pair (2 rows). One extra row (211) is appended at the bottom to hold
# the data that fell off the end of the block.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 97
$lastRowOld = 210
$lastRowNew = 211
$nOld = $lastRowOld - $firstRow + 1   # 114
$nNew = $lastRowNew - $firstRow + 1   # 115

# Snapshot all current values for the block (A:R) before writing anything.
$oldRange = $ws.Range("A$firstRow`:R$lastRowOld")
$old = $oldRange.Value2

# Prepare the output array (one extra row).
$new = New-Object 'object[,]' $nNew, 18

# Columns (1-based): A=1 B=2 C=3 D=4 E=5 F=6 G=7 H=8 I=9 J=10 K=11 L=12 M=13 N=14 O=15 P=16 Q=17 R=18
$qualityCols = 9, 10, 11, 12, 13, 15, 16   # I, J, K, L, M, O, P

for ($i = 1; $i -le $nNew; $i++) {
    $sheetRow = $firstRow + $i - 1

    if ($i -eq 1) {
        # Row 97: brand-new sample, copy the static columns from the old row 97
        # and set the new date/quality values explicitly.
        for ($c = 1; $c -le 18; $c++) { $new[$i - 1, $c - 1] = $old[$i, $c] }
        $new[$i - 1, 3]  = 44483   # D
        $new[$i - 1, 8]  = "Primera" # I
        $new[$i - 1, 9]  = 900     # J
        $new[$i - 1, 10] = 650     # K
        $new[$i - 1, 11] = 700     # L
        $new[$i - 1, 12] = 678     # M
        $new[$i - 1, 15] = 136     # P
    }
    else {
        # Base row: copy everything from the same old row first (static cols:
        # A,B,C,E,F,G,H,N,Q,R stay put), then overwrite the shifted columns.
        $baseIdx = [Math]::Min($i, $nOld)
        for ($c = 1; $c -le 18; $c++) { $new[$i - 1, $c - 1] = $old[$baseIdx, $c] }

        # Quality/volume/price/origin columns shift down by 1 row.
        foreach ($c in $qualityCols) {
            $new[$i - 1, $c - 1] = $old[$i - 1, $c]
        }

        # Fecha (D) is shared by pairs of rows; it shifts down by a full pair
        # (2 rows) on "odd" sheet rows and stays put on "even" sheet rows.
        if (($sheetRow % 2) -ne 0) {
            $new[$i - 1, 3] = $old[$i - 2, 4]
        }
    }
}

$destRange = $ws.Range("A$firstRow`:R$lastRowNew")
$destRange.Value2 = $new

# The brand-new row 211 needs the date-formatted style (matching column D in
# the rest of the block) applied to its D cell.
$ws.Range("D$lastRowNew").NumberFormat = $ws.Range("D$lastRowOld").NumberFormat
$ws.Range("D$lastRowNew").Value2 = $new[$nNew - 1, 3]
